$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 2-10 down into rows 3-11 (columns B:G),
# preserving column A (Q labels) and styles which are untouched.
# Iterate bottom-up so we don't clobber source rows before they are read.
for ($r = 10; $r -ge 2; $r--) {
    for ($c = 2; $c -le 7; $c++) {
        $v = $ws.Cells.Item($r, $c).Value2
        $ws.Cells.Item($r + 1, $c).Value2 = $v
    }
}

# New values for row 2 (top of the table)
$ws.Cells.Item(2, 2).Value2 = 0.1783908196033299
$ws.Cells.Item(2, 3).Value2 = 0.3606156554386025
$ws.Cells.Item(2, 4).Value2 = 0.2599511937740667
$ws.Cells.Item(2, 5).Value2 = 0.5098540906711122
$ws.Cells.Item(2, 6).Value2 = 0.4943913024279584
$ws.Cells.Item(2, 7).Value2 = 15
